# Updated Input AddFare -> now with decimal places
#
# 1) "edit" sheet: row 16's timestamp gets a tiny precision correction
#    (same instant, more accurate fractional-day representation).
# 2) "add" sheet: a new logged "AddFare" action row (row 8) is appended,
#    recording the fare Test01/Test10 with a 0.5 (decimal) value in column E.
# 3) "delete" sheet: a matching new row (row 10) is appended with the
#    combined description for that same action.

$wb = $excel.ActiveWorkbook

# --- 1) "edit" sheet -------------------------------------------------
$editWs = $wb.Worksheets.Item("edit")
$editWs.Range("A16").Value = 44413.70106576389

# --- 2) "add" sheet ----------------------------------------------------
$addWs = $wb.Worksheets.Item("add")

$addWs.Range("A8").Value = 44413.71606891204
$addWs.Range("A8").NumberFormat = "yyyy-mm-dd h:mm:ss"

$addWs.Range("B8").Value = "Test01"
$addWs.Range("C8").Value = "Test10"
$addWs.Range("D8").Value = 0
$addWs.Range("E8").Value = 0.5
$addWs.Range("F8").Value = 0
$addWs.Range("G8").Value = 0
$addWs.Range("H8").Value = 0
$addWs.Range("I8").Value = 0
$addWs.Range("J8").Value = 0
$addWs.Range("K8").Value = 0
$addWs.Range("L8").Value = 0

# --- 3) "delete" sheet --------------------------------------------------
$delWs = $wb.Worksheets.Item("delete")

$delWs.Range("A10").Value = 44413.71625334924
$delWs.Range("A10").NumberFormat = "yyyy-mm-dd h:mm:ss"

$delWs.Range("B10").Value = "Test01 Test10"
